$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.809.11"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.541.11"
$ws.Range("E3").Value = "  -1.73%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'206.12"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "'21.31"
$ws.Range("E9").Value = "  -2.89%  "
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").Value = "'0.0855"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").Value = "1.760.39"
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("D13").Value = "1.541.28"
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D16").Value = "26.811.25"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "'61.17"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "'214.03"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("E19").Value = "  -2.24%  "
$ws.Range("D20").Value = "0.0₃0682"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("E22").Value = "  -2.69%  "
$ws.Range("D23").Value = "'9.15"
$ws.Range("E23").Value = "  -1.77%  "
$ws.Range("E24").Value = "  -3.39%  "
$ws.Range("D25").Value = "'152.04"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").Value = "'6.59"
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("D27").Value = "'14.81"
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("D32").Value = "'3.23"
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("D33").Value = "1.368.28"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("E35").Value = "  -1.57%  "
$ws.Range("D36").Value = "'0.964"
$ws.Range("E36").Value = "  +3.50%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'0.805"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'5.75"
$ws.Range("E41").Value = "  +8.13%  "
$ws.Range("D42").Value = "'0.990"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").Value = "'2.21"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("D44").Value = "'63.00"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").Value = "'1.73"
$ws.Range("E45").Value = "  -3.63%  "
$ws.Range("D46").Value = "1.674.76"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").Value = "'84.31"
$ws.Range("E47").Value = "  -1.96%  "
$ws.Range("D49").Value = "0.0₇0977"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "'0.0942"
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("E51").Value = "  +0.02%  "
